$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark near the top of the document (Title
#    paragraph). Word re-numbers the remaining bookmarks automatically, which
#    also produces the "_Hlk70350700" (1->0) and "_Hlk70350350" (2->1)
#    re-numbering seen in the diff.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Fix the typo "thyroid-simulating" -> "thyroid-stimulating" (also clears
#    the spell-check proofErr markers that flagged the misspelling).
# ---------------------------------------------------------------------------
$fix = $d.Content
$fix.Find.Execute("thyroid-simulating hormone (TSH) testing every eight weeks, h", `
                   $true, $false, $false, $false, $false, $true, 1, $false, `
                   "thyroid-stimulating hormone (TSH) testing every eight weeks, h", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Re-insert a fresh "_GoBack" bookmark, landing (as Word does after the
#    user's last edit) in the middle of the word "Administered" inside the
#    "[\"Medication, Administered\": ...]" CQL snippet that follows the text
#    "we deleted it from the define ...".
# ---------------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("we deleted it from the define", $true, $false, $false, `
                      $false, $false, $true, 1, $false, "", 0) | Out-Null
$afterAnchor = $anchor.End

$target = $d.Range($afterAnchor, $d.Content.End)
$target.Find.Execute("Administered", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0) | Out-Null

$splitPoint = $target.Start + 3   # after "Adm", i.e. between "dm" and "inistered"
$bmRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
